$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns B, C, E keep their existing General/text style (they are naturally
# non-numeric so no special handling is required). Column D sometimes holds values
# that look like plain numbers (e.g. "1.012"); force those cells to Text format so
# Excel stores them as strings (matching the original inlineStr cells) instead of
# re-interpreting them as numbers, then restore the default "Normal" style so no
# extra style index lingers on the cell.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '29.601.75'
$ws.Range("E2").Value = '  +0.42%  '

# Row 3
$ws.Range("D3").Value = '1.930.58'
$ws.Range("E3").Value = '  +0.93%  '

# Row 4
$ws.Range("D4").Value = '1.012'
$ws.Range("E4").Value = '  +0.53%  '

# Row 5
$ws.Range("D5").Value = '326.49'
$ws.Range("E5").Value = '  +0.24%  '

# Row 6
$ws.Range("E6").Value = '  +0.44%  '

# Row 7
$ws.Range("D7").Value = '0.4829'
$ws.Range("E7").Value = '  +0.15%  '

# Row 8
$ws.Range("D8").Value = '0.4063'
$ws.Range("E8").Value = '  -0.04%  '

# Row 9
$ws.Range("D9").Value = '0.08227'
$ws.Range("E9").Value = '  +0.77%  '

# Row 10
$ws.Range("D10").Value = '1.014'
$ws.Range("E10").Value = '  -0.09%  '

# Row 11
$ws.Range("D11").Value = '23.90'
$ws.Range("E11").Value = '  +1.99%  '

# Row 12
$ws.Range("D12").Value = '1.961.05'
$ws.Range("E12").Value = '  +1.89%  '

# Row 13
$ws.Range("D13").Value = '6.124'
$ws.Range("E13").Value = '  +1.78%  '

# Row 14
$ws.Range("D14").Value = '7.332'
$ws.Range("E14").Value = '  +2.18%  '

# Row 15
$ws.Range("D15").Value = '91.75'
$ws.Range("E15").Value = '  +1.60%  '

# Row 16
$ws.Range("D16").Value = '0.06879'
$ws.Range("E16").Value = '  +1.21%  '

# Row 17
$ws.Range("D17").Value = '1.012'
$ws.Range("E17").Value = '  +0.40%  '

# Row 18
$ws.Range("E18").Value = '  +0.14%  '

# Row 19
$ws.Range("D19").Value = '17.70'
$ws.Range("E19").Value = '  -0.06%  '

# Row 20
$ws.Range("E20").Value = '  +0.44%  '

# Row 21
$ws.Range("D21").Value = '29.593.44'
$ws.Range("E21").Value = '  +0.29%  '

# Row 22
$ws.Range("D22").Value = '5.687'
$ws.Range("E22").Value = '  +0.96%  '

# Row 23
$ws.Range("D23").Value = '12.07'
$ws.Range("E23").Value = '  +2.82%  '

# Row 24
$ws.Range("D24").Value = '2.184'
$ws.Range("E24").Value = '  +0.08%  '

# Row 25
$ws.Range("D25").Value = '2.173.28'
$ws.Range("E25").Value = '  +1.02%  '

# Row 26
$ws.Range("D26").Value = '156.07'
$ws.Range("E26").Value = '  +0.22%  '

# Row 27
$ws.Range("D27").Value = '6.423'
$ws.Range("E27").Value = '  -0.35%  '

# Row 28
$ws.Range("D28").Value = '20.08'
$ws.Range("E28").Value = '  +0.13%  '

# Row 29
$ws.Range("D29").Value = '2.105'
$ws.Range("E29").Value = '  -0.38%  '

# Row 30
$ws.Range("D30").Value = '121.06'
$ws.Range("E30").Value = '  +0.71%  '

# Row 31
$ws.Range("D31").Value = '1.016'
$ws.Range("E31").Value = '  -0.91%  '

# Row 32
$ws.Range("D32").Value = '0.09610'
$ws.Range("E32").Value = '  +0.73%  '

# Row 33
$ws.Range("D33").Value = '5.610'
$ws.Range("E33").Value = '  +1.37%  '

# Row 34
$ws.Range("D34").Value = '3.560'
$ws.Range("E34").Value = '  -0.14%  '

# Row 35
$ws.Range("D35").Value = '1.387'
$ws.Range("E35").Value = '  -0.35%  '

# Row 36
$ws.Range("D36").Value = '0.06397'
$ws.Range("E36").Value = '  +4.83%  '

# Row 37
$ws.Range("D37").Value = '0.02292'
$ws.Range("E37").Value = '  +0.83%  '

# Row 38
$ws.Range("D38").Value = '1.191'
$ws.Range("E38").Value = '  +1.24%  '

# Row 39
$ws.Range("D39").Value = '0.5971'
$ws.Range("E39").Value = '  +0.08%  '

# Row 40
$ws.Range("E40").Value = '  -0.99%  '

# Row 41
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '7.915'
$ws.Range("E41").Value = '  -1.22%  '

# Row 42
$ws.Range("B42").Value = 'Frax'
$ws.Range("C42").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D42").Value = '1.011'
$ws.Range("E42").Value = '  +0.42%  '

# Row 43
$ws.Range("D43").Value = '0.1851'
$ws.Range("E43").Value = '  -0.21%  '

# Row 44
$ws.Range("D44").Value = '2.488'
$ws.Range("E44").Value = '  +4.24%  '

# Row 45
$ws.Range("D45").Value = '1.247'
$ws.Range("E45").Value = '  -2.55%  '

# Row 46
$ws.Range("D46").Value = '12.43'
$ws.Range("E46").Value = '  -0.97%  '

# Row 47
$ws.Range("D47").Value = '0.07537'
$ws.Range("E47").Value = '  -0.91%  '

# Row 48
$ws.Range("D48").Value = '0.5572'
$ws.Range("E48").Value = '  -0.14%  '

# Row 49
$ws.Range("E49").Value = '  +1.71%  '

# Row 50
$ws.Range("D50").Value = '119.34'
$ws.Range("E50").Value = '  +2.69%  '

# Row 51
$ws.Range("D51").Value = '2.441'
$ws.Range("E51").Value = '  +1.56%  '

# Restore default style on column D so no stray number-format style remains
$dRange.Style = "Normal"
